$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("H2").Value = 92
$ws.Range("I2").Value = 245
$ws.Range("J2").Value = 922
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 269
$ws.Range("M2").Value = 17
$ws.Range("N2").Value = 166
$ws.Range("P2").Value = 2
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 11
$ws.Range("T2").Value = 187
$ws.Range("U2").Value = 11
$ws.Range("V2").Value = 1429
$ws.Range("W2").Value = 0
$ws.Range("X2").Value = 1504
$ws.Range("Y2").Value = 1
$ws.Range("Z2").Value = 22
$ws.Range("AA2").Value = 10
